$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of fixtures on this sheet had their "payload" (column B id plus
# the odds/result columns F:AC) re-associated with the correct row, while
# columns A (rank), C, D, E (league/date) stay put on their original row.
#
# Rows 192 and 193 simply swap their B,F:AC payload (2-cycle).
# Rows 204,205,206,207 rotate their B,F:AC payload cyclically:
#   new(204) = old(205); new(205) = old(206); new(206) = old(207); new(207) = old(204)
#
# NOTE: avoid helper functions here - returning a 2D Range.Value() array out
# of a PowerShell function mangles/flattens it in this runtime, so everything
# is done inline with plain variables instead.

# --- Swap rows 192 and 193 ---
$b192 = $ws.Range("B192").Value()
$fac192 = $ws.Range("F192:AC192").Value()
$b193 = $ws.Range("B193").Value()
$fac193 = $ws.Range("F193:AC193").Value()

$ws.Range("B192").Value = $b193
$ws.Range("F192:AC192").Value = $fac193
$ws.Range("B193").Value = $b192
$ws.Range("F193:AC193").Value = $fac192

# --- Rotate rows 204 -> 205 -> 206 -> 207 -> 204 ---
$b204 = $ws.Range("B204").Value()
$fac204 = $ws.Range("F204:AC204").Value()
$b205 = $ws.Range("B205").Value()
$fac205 = $ws.Range("F205:AC205").Value()
$b206 = $ws.Range("B206").Value()
$fac206 = $ws.Range("F206:AC206").Value()
$b207 = $ws.Range("B207").Value()
$fac207 = $ws.Range("F207:AC207").Value()

# new(204) = old(205)
$ws.Range("B204").Value = $b205
$ws.Range("F204:AC204").Value = $fac205

# new(205) = old(206)
$ws.Range("B205").Value = $b206
$ws.Range("F205:AC205").Value = $fac206

# new(206) = old(207)
$ws.Range("B206").Value = $b207
$ws.Range("F206:AC206").Value = $fac207

# new(207) = old(204)
$ws.Range("B207").Value = $b204
$ws.Range("F207:AC207").Value = $fac204
